$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "month of plan" dates in column A are stored as plain text
# (e.g. "2008-01-01"), not real Excel dates. Force the cells to stay
# text (so Excel doesn't silently reinterpret them as date serials)
# while updating the values, then clear the temporary formatting so
# the cells keep their original (unstyled) appearance.
$rng = $ws.Range("A2:A7")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "2009-01-01"
$ws.Range("A3").Value = "2026-03-01"
$ws.Range("A4").Value = "2026-04-01"
$ws.Range("A5").Value = "2027-04-01"
$ws.Range("A6").Value = "2028-05-01"
$ws.Range("A7").Value = "2029-06-01"

$rng.ClearFormats()
